# Update the date heading
$d = $word.ActiveDocument
$d.Content.Find.Execute("2025-01-30 Thursday", $true, $false, $false, $false, $false, $true, 1, $false, "2025-01-31 Friday", 2)

# Update the answer table (first table in the document)
$t = $d.Tables(1)

# Row 1: drop "53÷5=10, 3" (cells shift left), update remaining cells,
# and append a new trailing cell "30÷7=4, 2"
$t.Cell(1, 1).Range.Text = "16÷3=5, 1"
$t.Cell(1, 2).Range.Text = "88÷8=11, 0"
$t.Cell(1, 3).Range.Text = "44÷6=7, 2"
$t.Cell(1, 4).Range.Text = "44÷7=6, 2"
$t.Cell(1, 5).Range.Text = "30÷7=4, 2"

# Row 5
$t.Cell(5, 1).Range.Text = "76÷9=8, 4"
$t.Cell(5, 2).Range.Text = "97÷7=13, 6"
$t.Cell(5, 3).Range.Text = "52÷6=8, 4"
$t.Cell(5, 4).Range.Text = "92÷6=15, 2"
$t.Cell(5, 5).Range.Text = "16÷5=3, 1"

# Row 9
$t.Cell(9, 1).Range.Text = "32÷7=4, 4"
$t.Cell(9, 2).Range.Text = "56÷8=7, 0"
$t.Cell(9, 3).Range.Text = "96÷6=16, 0"
$t.Cell(9, 4).Range.Text = "10÷6=1, 4"
$t.Cell(9, 5).Range.Text = "22÷4=5, 2"

# Row 13
$t.Cell(13, 1).Range.Text = "15÷7=2, 1"
$t.Cell(13, 2).Range.Text = "16÷8=2, 0"
$t.Cell(13, 3).Range.Text = "82÷9=9, 1"
$t.Cell(13, 4).Range.Text = "42÷9=4, 6"
$t.Cell(13, 5).Range.Text = "82÷6=13, 4"

# Row 17
$t.Cell(17, 1).Range.Text = "78÷3=26, 0"
$t.Cell(17, 2).Range.Text = "78÷5=15, 3"
$t.Cell(17, 3).Range.Text = "89÷7=12, 5"
$t.Cell(17, 4).Range.Text = "89÷7=12, 5"
$t.Cell(17, 5).Range.Text = "80÷8=10, 0"
